$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a common "fight" event to the QuestDungeon (column M) lists for each
# dungeon row, matching commit "add common fight event to dungeon".
# Values are set in the same order the original author touched them so the
# shared-string table gets rebuilt with matching new-entry order.
$ws.Range("M5").Value = "fight;5|trees;4"
$ws.Range("M4").Value = "fight;7|trees;2|manflower;2|river;2|cliff;2|losttree;1|oldtree;1"
$ws.Range("M6").Value = "fight;10|sandland;2|potteryroom;2|honeyhome;2|snare;1|basement;1|woodhouse2;1|booty;1|trapspear;2|trapdrop;1|potteryman;1|stonedoor2;1|crystalball;2"
$ws.Range("M7").Value = "fight;5|trees;4"

# Move the active selection to M6, matching the saved cursor position in the
# edited workbook.
$ws.Range("M6").Select()

$wb.Save()
